$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bug fix: the Monday hours for the week commencing 05/03/2018 (row 8) were
# missing from the sheet. Add the missing value; the Total column (I8) and
# the grand total (I19) recalculate automatically via their formulas.
$ws.Range("B8").Value = 4.25

# Move the active selection to F12, matching where the user left off editing.
$ws.Range("F12").Select()
